$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medium traffic density")

# --- Block 1: rows 3-13 ---
$ws.Range("D3").Value = 3792.94
$ws.Range("E3").Value = 311.7

$ws.Range("D4").Value = 3544.49
$ws.Range("E4").Value = 300.45

$ws.Range("D8").Value = 3674.95
$ws.Range("E8").Value = 324.43

$ws.Range("D9").Value = 3622.18
$ws.Range("E9").Value = 300.22

$ws.Range("D10").Value = 3624.6
$ws.Range("E10").Value = 301.559

$ws.Range("D11").Value = 3650.35
$ws.Range("E11").Value = 310.95

$ws.Range("D12").Value = 3524.76
$ws.Range("E12").Value = 300.74

$ws.Range("D13").Formula = "=AVERAGE(D8:D12)"
$ws.Range("E13").Formula = "=AVERAGE(E8:E12)"

# --- Block 2: rows 17-27 ---
$ws.Range("D17").Value = 1868.03
$ws.Range("E17").Value = 233.11

$ws.Range("D22").Value = 1783.16
$ws.Range("E22").Value = 233.719

$ws.Range("D23").Value = 1725.47
$ws.Range("E23").Value = 238.136

$ws.Range("D24").Value = 1720.94
$ws.Range("E24").Value = 235.07

$ws.Range("D25").Value = 1658.64
$ws.Range("E25").Value = 230.76

$ws.Range("D26").Value = 1659
$ws.Range("E26").Value = 222.508

$ws.Range("D27").Formula = "=AVERAGE(D22:D26)"
$ws.Range("E27").Formula = "=AVERAGE(E22:E26)"

# --- Selection change ---
$ws.Range("D23").Select()
